$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Backtick character used to build the embedded Neo4j-style `backtick` identifiers
$bt = [char]96

$lines = @(
    "MATCH (p:participant)-->(s:study)",
    "OPTIONAL MATCH (samp:sample)-->(p)",
    "OPTIONAL MATCH (p)<--(diag:diagnosis)",
    "OPTIONAL MATCH (samp)<--(f:file)",
    "OPTIONAL MATCH (f)<--(g:genomic_info)",
    "WITH s, p, samp, f, g, diag",
    "WHERE f.file_type in ['BAI']",
    "with p",
    "OPTIONAL MATCH (p)-->(s:study)",
    "OPTIONAL MATCH (samp:sample)-->(p)",
    "WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp",
    "RETURN",
    "coalesce(p.participant_id,'') as ${bt}Participant ID${bt},",
    "coalesce(s.study_name, '') as ${bt}Study Name${bt},",
    "coalesce(s.phs_accession,'') as ${bt}Accession${bt},",
    "coalesce(p.gender,'') as ${bt}Gender${bt},",
    "coalesce(apoc.text.join(samp, ','), '') as ${bt}Samples${bt}",
    "ORDER BY p.participant_id LIMIT 100"
)
$newQuery = [string]::Join("`n", $lines)

# Update the query text stored in cell B2 (the "ParticipantsTab" query)
$ws.Range("B2").Value2 = $newQuery

# The longer query text wraps across more lines, so the row grows taller
$ws.Rows.Item(2).RowHeight = 279
